$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 5 Hydro-Met sub-program rows that were removed:
# "National Water Data Program", "Drought Monitoring and Impact Assessment",
# "Flood Forecasting and Warning ", "Agrometeorological Advisory Services",
# "National Climate Assessment " (rows 26-30)
$ws.Rows("26:30").Delete()

$ws.Range("B36").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
